$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 972.5
$ws.Range("I19").Value = 973
$ws.Range("J19").Value = 972
$ws.Range("K19").Value = 973
$ws.Range("L19").Value = 972
$ws.Range("M19").Value = -798
$ws.Range("N19").Value = -1322
$ws.Range("H43").Value = 1692
$ws.Range("I43").Value = 1692
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1692
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1623
$ws.Range("N43").ClearContents()
$ws.Range("H69").Value = 34917.668
$ws.Range("J69").Value = 18090.818
$ws.Range("L69").Value = 54272.454
$ws.Range("N69").Value = -56020.454
$ws.Range("H72").Value = 34917.668
$ws.Range("J72").Value = 18090.818
$ws.Range("L72").Value = 162817.362
$ws.Range("N72").Value = -171553.362
$ws.Range("H80").Value = 8244.6
$ws.Range("I80").Value = 5300.6665
$ws.Range("K80").Value = 15901.9995
$ws.Range("M80").Value = -14903.9995
$ws.Range("H83").Value = 8244.6
$ws.Range("I83").Value = 5300.6665
$ws.Range("K83").Value = 47705.9985
$ws.Range("M83").Value = -42713.9985
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 1000
$ws.Range("K94").Value = 1000
$ws.Range("M94").Value = -549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 579.8182
$ws.Range("I2").Value = 387.8
$ws.Range("K2").Value = 387.8
$ws.Range("M2").Value = -274.8
$ws.Range("H46").Value = 10045.667
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H110").Value = 1328.1666
$ws.Range("I110").Value = 1328.1666
$ws.Range("K110").Value = 1328.1666
$ws.Range("M110").Value = 716.8334
$ws.Range("H116").Value = 579.8182
$ws.Range("I116").Value = 387.8
$ws.Range("K116").Value = 387.8
$ws.Range("M116").Value = 1906.2
$ws.Range("H132").Value = 1974.7273
$ws.Range("I132").Value = 1974.7273
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5924.1819
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3394.1819
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 579.8182
$ws.Range("I3").Value = 387.8
$ws.Range("K3").Value = 387.8
$ws.Range("M3").Value = -273.8
$ws.Range("H22").Value = 677.5
$ws.Range("I22").Value = 540
$ws.Range("J22").Value = 1640
$ws.Range("K22").Value = 540
$ws.Range("L22").Value = 1640
$ws.Range("M22").Value = -367
$ws.Range("N22").Value = -1986
$ws.Range("H81").Value = 40755
$ws.Range("J81").Value = 40755
$ws.Range("L81").Value = 40755
$ws.Range("N81").Value = -42877
$ws.Range("H84").Value = 40755
$ws.Range("J84").Value = 40755
$ws.Range("L84").Value = 122265
$ws.Range("N84").Value = -132873
$ws.Range("H86").Value = 2268
$ws.Range("I86").Value = 2268
$ws.Range("K86").Value = 2268
$ws.Range("M86").Value = -1145
$ws.Range("H89").Value = 2268
$ws.Range("I89").Value = 2268
$ws.Range("K89").Value = 11340
$ws.Range("M89").Value = -5724
$ws.Range("H99").Value = 2411.077
$ws.Range("I99").Value = 1893.3
$ws.Range("J99").Value = 4137
$ws.Range("K99").Value = 1893.3
$ws.Range("L99").Value = 4137
$ws.Range("M99").Value = -395.3
$ws.Range("N99").Value = -7133

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2978.111
$ws.Range("I16").Value = 3405.5
$ws.Range("J16").Value = 2636.2
$ws.Range("K16").Value = 3405.5
$ws.Range("L16").Value = 2636.2
$ws.Range("M16").Value = -3118.5
$ws.Range("N16").Value = -3210.2
$ws.Range("H35").Value = 1600
$ws.Range("I35").Value = 1733.3334
$ws.Range("K35").Value = 1733.3334
$ws.Range("M35").Value = -1439.3334
$ws.Range("H113").Value = 2978.111
$ws.Range("I113").Value = 3405.5
$ws.Range("J113").Value = 2636.2
$ws.Range("K113").Value = 3405.5
$ws.Range("L113").Value = 2636.2
$ws.Range("M113").Value = -1235.5
$ws.Range("N113").Value = -6976.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1294.75
$ws.Range("I9").Value = 180
$ws.Range("J9").Value = 1666.3334
$ws.Range("K9").Value = 540
$ws.Range("L9").Value = 4999.0002
$ws.Range("M9").Value = -316
$ws.Range("N9").Value = -5447.0002
$ws.Range("H119").Value = 3985
$ws.Range("I119").Value = 3985
$ws.Range("K119").Value = 11955
$ws.Range("M119").Value = -7117
$ws.Range("H121").Value = 14317
$ws.Range("J121").Value = 7129.3335
$ws.Range("L121").Value = 21388.0005
$ws.Range("N121").Value = -24008.0005
$ws.Range("H131").Value = 3186.625
$ws.Range("J131").Value = 3186.625
$ws.Range("L131").Value = 9559.875
$ws.Range("N131").Value = -19639.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1000
$ws.Range("J80").Value = 1000
$ws.Range("L80").Value = 1000
$ws.Range("N80").Value = -2996
$ws.Range("H83").Value = 1000
$ws.Range("J83").Value = 1000
$ws.Range("L83").Value = 5000
$ws.Range("N83").Value = -14984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 408.45
$ws.Range("I55").Value = 363
$ws.Range("J55").Value = 666
$ws.Range("K55").Value = 363
$ws.Range("L55").Value = 666
$ws.Range("M55").Value = -190
$ws.Range("N55").Value = -1012
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 2942.6667
$ws.Range("I132").Value = 2481.611
$ws.Range("K132").Value = 7444.833
$ws.Range("M132").Value = -4914.833
$ws.Range("H136").Value = 3053.261
$ws.Range("I136").Value = 2751.75
$ws.Range("K136").Value = 8255.25
$ws.Range("M136").Value = -5705.25
$ws.Range("H137").Value = 110386.5
$ws.Range("I137").Value = 110386
$ws.Range("J137").Value = 110387
$ws.Range("K137").Value = 110386
$ws.Range("L137").Value = 110387
$ws.Range("M137").Value = -105286
$ws.Range("N137").Value = -120587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14999
$ws.Range("J62").Value = 14999
$ws.Range("L62").Value = 14999
$ws.Range("N62").Value = -16247
$ws.Range("H65").Value = 14999
$ws.Range("J65").Value = 14999
$ws.Range("L65").Value = 74995
$ws.Range("N65").Value = -81235
$ws.Range("H126").Value = 4710.6313
$ws.Range("I126").Value = 3214.7144
$ws.Range("K126").Value = 9644.143199999999
$ws.Range("M126").Value = -7174.143199999999
$ws.Range("H136").Value = 2818.7273
$ws.Range("I136").Value = 2667.3333
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 8001.999899999999
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -5451.999899999999
